$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("tests")

# Remove the old highlighted row 6 ("buck filtering" entry with extra G/H cells,
# special yellow fill + custom number format). This shifts rows 7-8 up to 6-7.
$ws.Rows.Item(6).Delete()

# Row 7 (was old row 8) keeps its data but the test date in column C changes.
$ws.Range("C7").Value = 41044

# Re-add a "buck filtering" entry as the new last row (row 8) -- the
# impedance analyzer files entry -- without the old yellow highlighting.
$ws.Range("A8").Value = 41045
$ws.Range("B8").Value = "buck filtering"
$ws.Range("C8").Value = 41045
$ws.Range("D8").Value = "The capacitors were put in a low power buck supply to test for degredation"
$ws.Range("E8").Value = "0x0002"
$ws.Range("F8").Value = "0x0003"

# Match the plain date style already used by the other date cells (A3, C3, ...).
$ws.Range("A3").Copy()
$ws.Range("A8").PasteSpecial(-4122)
$ws.Range("C3").Copy()
$ws.Range("C8").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("C7").Select()
